$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# Date updated
$ws1.Range("B8").Value = "2023-09-01T18:07:06+00:00"

# Content: supplement -> complete
$ws1.Range("B19").Value = "complete"

# Count: (empty) -> 1  (copy from Concepts!A2 which already holds text "1",
# so the destination keeps its original style and the cell is written as a
# shared string rather than a number)
$ws2.Range("A2").Copy($ws1.Range("B21"))
